$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update existing row 2 values
$ws.Range("B2").Value = "2023-10-16 18:06:27"
$ws.Range("C2").Value = "172.29.0.1"

# Add new row 3
$ws.Range("A3").Value = "Bartek"
$ws.Range("B3").Value = "2023-10-16 18:16:00"
$ws.Range("C3").Value = "192.168.16.1"
